$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C130").Value = "[name=`"Talulah`"]  The 'Emperors' Blade'...? *Spit*! You came all the way here today to mock me? Or are you here to kill me?!   `n"
$ws.Range("C137").Value = "[name=`"'Emperors' Blade'`"]  *Hiss*.`n"
$ws.Range("C138").Value = "[name=`"'Emperors' Blade'`"]  You didn't tell them who you really are.`n"
$ws.Range("C139").Value = "[name=`"'Emperors' Blade'`"]  Is this part of your plan?`n"
$ws.Range("C143").Value = "[name=`"'Emperors' Blade'`"]  You can only turn one against another if there existed trust between them in the first place. I doubt there is any trust between you and the rest of them. `n"
$ws.Range("C144").Value = "[name=`"'Emperors' Blade'`"]  I can only surmise... You believe the trust they have for you will remain unscathed even after they come to know who you really are.`n"
$ws.Range("C146").Value = "[name=`"'Emperors' Blade'`"]  Your words are strong... and confident.`n"
$ws.Range("C147").Value = "[name=`"'Emperors' Blade'`"]  If this is as you expected... Very well, I must reassess the situation.`n"
$ws.Range("C148").Value = "[name=`"'Emperors' Blade'`"]  I bid you goodbye, then. Remember this well. There is one thing we share with the people standing around you. `n"
$ws.Range("C149").Value = "[name=`"'Emperors' Blade'`"]  We have no trust in you. You must take action to earn even the tiniest trickle of trust. `n"
$ws.Range("C152").Value = "[name=`"'Emperors' Blade'`"]  *Hiss*... You have a question?`n"
$ws.Range("C154").Value = "[name=`"'Emperors' Blade'`"]  Am I perhaps misunderstanding you?`n"
$ws.Range("C156").Value = "[name=`"'Emperors' Blade'`"]  *Huff*...`n"
$ws.Range("C157").Value = "[name=`"'Emperors' Blade'`"]  Which secret?`n"
$ws.Range("C166").Value = "[name=`"'Emperors' Blade'`"]  You are mistaken about one thing... I am not here alone.`n"
$ws.Range("C167").Value = "[name=`"'Emperors' Blade'`"]  There is one more of us behind you.`n"
$ws.Range("C173").Value = "[name=`"'Emperors' Blade'`"]  Daughter of Kashchey, allow me to postulate that you've come to the following conclusion: Your peers will not be suspicious of you merely because of your identity.`n"
$ws.Range("C174").Value = "[name=`"'Emperors' Blade'`"]  Then let us verify the results. If your judgment proves incorrect, then the promise your father made to us shall be considered void. `n"
$ws.Range("C187").Value = "[name=`"'Emperors' Blade'`"]  What you said could not be more different than how your father would see it.`n"
$ws.Range("C188").Value = "[name=`"'Emperors' Blade'`"]  ...I am afraid you may only prove to be a disappointment, leader of the Northern Tundra Infected.`n"
